# Updates cryptos list values (Price/Volume columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "4.439",
    # "0.5213", "26.133.26") are not auto-coerced into numbers/dates,
    # matching the original inline-string cell type.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    # Restore default (un-styled) cell formatting so the style index
    # matches the source workbook (these cells carry no explicit "s").
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.133.26"
Set-TextValue "D3" "1.669.69"
Set-TextValue "E3" "  -1.42%  "
Set-TextValue "E4" "  -0.74%  "
Set-TextValue "D5" "210.62"
Set-TextValue "E5" "  -3.93%  "
Set-TextValue "D6" "0.5213"
Set-TextValue "E6" "  -5.00%  "
Set-TextValue "E7" "  -0.75%  "
Set-TextValue "D8" "0.2638"
Set-TextValue "E8" "  -3.53%  "
Set-TextValue "D9" "0.06247"
Set-TextValue "E9" "  -3.48%  "
Set-TextValue "E10" "  -4.05%  "
Set-TextValue "D11" "0.07509"
Set-TextValue "E11" "  -2.21%  "
Set-TextValue "D12" "1.693.20"
Set-TextValue "E12" "  -0.13%  "
Set-TextValue "D13" "4.439"
Set-TextValue "E14" "  -4.56%  "
Set-TextValue "D15" "0.000007971"
Set-TextValue "E15" "  -4.98%  "
Set-TextValue "D16" "66.21"
Set-TextValue "E16" "  +1.11%  "
Set-TextValue "D17" "26.183.04"
Set-TextValue "E17" "  -1.03%  "
Set-TextValue "E18" "  -0.75%  "
Set-TextValue "D19" "4.788"
Set-TextValue "E19" "  -3.34%  "
Set-TextValue "D20" "187.07"
Set-TextValue "E20" "  -2.86%  "
Set-TextValue "E21" "  -5.68%  "
Set-TextValue "D22" "6.174"
Set-TextValue "D23" "1.004"
Set-TextValue "E23" "  -0.71%  "
Set-TextValue "D24" "147.84"
Set-TextValue "E24" "  -0.99%  "
Set-TextValue "D25" "0.1242"
Set-TextValue "E25" "  -6.60%  "
Set-TextValue "D26" "7.582"
Set-TextValue "E26" "  -4.24%  "
Set-TextValue "D27" "15.90"
Set-TextValue "E27" "  +0.51%  "
Set-TextValue "D28" "0.06183"
Set-TextValue "E28" "  -1.93%  "
Set-TextValue "D29" "1.355"
Set-TextValue "E29" "  -3.00%  "
Set-TextValue "E30" "  -3.98%  "
Set-TextValue "D31" "3.472"
Set-TextValue "E31" "  -3.89%  "
Set-TextValue "E32" "  -5.02%  "
Set-TextValue "D33" "1.610"
Set-TextValue "E33" "  -4.49%  "
Set-TextValue "D34" "0.9909"
Set-TextValue "E34" "  -5.18%  "
Set-TextValue "D35" "0.6037"
Set-TextValue "E35" "  -1.91%  "
Set-TextValue "D36" "2.404"
Set-TextValue "E36" "  -0.19%  "
Set-TextValue "E37" "  -0.14%  "
Set-TextValue "E38" "  -1.15%  "
Set-TextValue "D39" "0.01608"
Set-TextValue "E39" "  -1.89%  "
Set-TextValue "D40" "1.069.82"
Set-TextValue "E40" "  -4.46%  "
Set-TextValue "D41" "0.8645"
Set-TextValue "E41" "  -2.48%  "
Set-TextValue "E42" "  -1.17%  "
Set-TextValue "E43" "  -2.44%  "
Set-TextValue "D44" "1.819.16"
Set-TextValue "E44" "  -1.37%  "
Set-TextValue "D45" "0.00000000110"
Set-TextValue "E45" "  +0.03%  "
Set-TextValue "D46" "55.94"
Set-TextValue "E47" "  -0.46%  "
Set-TextValue "D48" "0.05246"
Set-TextValue "E48" "  -0.80%  "
Set-TextValue "D49" "7.942"
Set-TextValue "E49" "  -3.52%  "
Set-TextValue "D50" "0.4251"
Set-TextValue "E50" "  -1.25%  "
Set-TextValue "D51" "5.945"
Set-TextValue "E51" "  -2.69%  "
